{"js": "// Replace the \"Summary 2\" table header text with \"Summary 0\".\nconst results = context.document.body.search(\"Summary 2\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"Summary 0\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$find = $d.Content.Find\n$find.Execute(\"Summary 2\", $false, $false, $false, $false, $false, $true, 1, $false, \"Summary 0\", 2)\n"}
